# Generate Report for Handback
# Refresh the handback-status report timestamps/priority after a re-run of
# the localization handback generation for the 004b7a50 (.md) and
# 6d4e8c76 (.md) source files.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" for the two files that
#     were just regenerated (rows 2 and 4: 004b7a50-...md and 6d4e8c76-...md)
$wsOverview.Range("G2").Value = "2016-08-21 00:15:56"
$wsOverview.Range("G4").Value = "2016-08-21 00:15:56"

# --- zh-cn sheet: Priority moved from human translation (ht) to machine
#     translation (mt), and the handoff/handback timestamps advanced.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

$wsZhCn.Range("H2").Value = "2016-08-21 00:15:52"
$wsZhCn.Range("H4").Value = "2016-08-21 00:15:52"

$wsZhCn.Range("K2").Value = "2016-08-21 00:16:13"
$wsZhCn.Range("K4").Value = "2016-08-21 00:16:13"

# --- de-de sheet: same priority change, plus correspond handoff/handback
#     datetime updates.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

$wsDeDe.Range("H2").Value = "2016-08-21 00:15:56"
$wsDeDe.Range("H4").Value = "2016-08-21 00:15:56"

$wsDeDe.Range("K2").Value = "2016-08-21 00:16:20"
$wsDeDe.Range("K4").Value = "2016-08-21 00:16:20"
